$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New section title
$ws.Range("A9").Value = "Comparing Spiral vs Structured using sum(1/d) criteria"

# New header row
$ws.Range("A11").Value = "N"
$ws.Range("B11").Value = "Spiral"
$ws.Range("C11").Value = "Structured"

# New data rows
$ws.Range("A12").Value = 4
$ws.Range("B12").Value = 2.6423254318158702
$ws.Range("C12").Value = 2.7316634235818902

$ws.Range("A13").Value = 9
$ws.Range("B13").Value = 8.8148050166905101
$ws.Range("C13").Value = 8.3793569722851799

$ws.Range("A14").Value = 16
$ws.Range("B14").Value = 19.1080983555775
$ws.Range("C14").Value = 16.971715494995902

$ws.Range("A15").Value = 25
$ws.Range("B15").Value = 33.933145084621202
$ws.Range("C15").Value = 28.508849662534701

# Apply the existing "0.0" number format style to the new numeric range (style index 1 in target)
$ws.Range("B12:C15").NumberFormat = "0.0"

# Update selection to match target sheetView
$ws.Range("A11:C15").Select()
